$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "partial"/"dissolution" columns (C, D) and old header/data values ---
$ws.Range("C1:D8").Clear()

# --- New headers ---
$ws.Range("A1").Value = "CTD"
$ws.Range("B1").Value = "Outside carapace"

# --- New station / dissolution-fraction data ---
$data = @(
  @(94, 0.5),
  @(90, 0.4),
  @(99, 0.75),
  @(100, 0.75),
  @(106, 0.4),
  @(109, 0.5),
  @(115, 0.1),
  @(128, 0)
)
$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $r++
}

# --- Header formatting: bold across the header row (A1:C1 and F1:H1) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("F1:H1").Font.Bold = $true

# --- Highlight the station column (A) ---
$ws.Range("A2:A9").Interior.Color = 65535

# --- Stray formatting on F1 (white background) ---
$ws.Range("F1").Interior.Color = 16777215

# --- Column widths ---
$ws.Columns(2).ColumnWidth = 12.71
$ws.Columns(3).ColumnWidth = 11.71

# --- Selection / view state ---
$ws.Range("A1:B9").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
